$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Formula = "=""28.361.47"""
$ws.Cells.Item(2, 4).Copy()
$ws.Cells.Item(2, 4).PasteSpecial(-4163)
$ws.Cells.Item(2, 5).Formula = "=""  +5.30%  """
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(2, 5).PasteSpecial(-4163)
$ws.Cells.Item(3, 4).Formula = "=""1.811.18"""
$ws.Cells.Item(3, 4).Copy()
$ws.Cells.Item(3, 4).PasteSpecial(-4163)
$ws.Cells.Item(3, 5).Formula = "=""  +4.66%  """
$ws.Cells.Item(3, 5).Copy()
$ws.Cells.Item(3, 5).PasteSpecial(-4163)
$ws.Cells.Item(4, 4).Formula = "=""1.001"""
$ws.Cells.Item(4, 4).Copy()
$ws.Cells.Item(4, 4).PasteSpecial(-4163)
$ws.Cells.Item(4, 5).Formula = "=""  +0.27%  """
$ws.Cells.Item(4, 5).Copy()
$ws.Cells.Item(4, 5).PasteSpecial(-4163)
$ws.Cells.Item(5, 4).Formula = "=""318.46"""
$ws.Cells.Item(5, 4).Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4163)
$ws.Cells.Item(5, 5).Formula = "=""  +2.83%  """
$ws.Cells.Item(5, 5).Copy()
$ws.Cells.Item(5, 5).PasteSpecial(-4163)
$ws.Cells.Item(6, 4).Formula = "=""1.001"""
$ws.Cells.Item(6, 4).Copy()
$ws.Cells.Item(6, 4).PasteSpecial(-4163)
$ws.Cells.Item(6, 5).Formula = "=""  +0.18%  """
$ws.Cells.Item(6, 5).Copy()
$ws.Cells.Item(6, 5).PasteSpecial(-4163)
$ws.Cells.Item(7, 4).Formula = "=""0.5704"""
$ws.Cells.Item(7, 4).Copy()
$ws.Cells.Item(7, 4).PasteSpecial(-4163)
$ws.Cells.Item(7, 5).Formula = "=""  +16.15%  """
$ws.Cells.Item(7, 5).Copy()
$ws.Cells.Item(7, 5).PasteSpecial(-4163)
$ws.Cells.Item(8, 4).Formula = "=""0.3875"""
$ws.Cells.Item(8, 4).Copy()
$ws.Cells.Item(8, 4).PasteSpecial(-4163)
$ws.Cells.Item(8, 5).Formula = "=""  +10.18%  """
$ws.Cells.Item(8, 5).Copy()
$ws.Cells.Item(8, 5).PasteSpecial(-4163)
$ws.Cells.Item(9, 4).Formula = "=""43.17"""
$ws.Cells.Item(9, 4).Copy()
$ws.Cells.Item(9, 4).PasteSpecial(-4163)
$ws.Cells.Item(9, 5).Formula = "=""  +0.61%  """
$ws.Cells.Item(9, 5).Copy()
$ws.Cells.Item(9, 5).PasteSpecial(-4163)
$ws.Cells.Item(10, 4).Formula = "=""0.07609"""
$ws.Cells.Item(10, 4).Copy()
$ws.Cells.Item(10, 4).PasteSpecial(-4163)
$ws.Cells.Item(10, 5).Formula = "=""  +4.92%  """
$ws.Cells.Item(10, 5).Copy()
$ws.Cells.Item(10, 5).PasteSpecial(-4163)
$ws.Cells.Item(11, 4).Formula = "=""1.139"""
$ws.Cells.Item(11, 4).Copy()
$ws.Cells.Item(11, 4).PasteSpecial(-4163)
$ws.Cells.Item(11, 5).Formula = "=""  +8.21%  """
$ws.Cells.Item(11, 5).Copy()
$ws.Cells.Item(11, 5).PasteSpecial(-4163)
$ws.Cells.Item(12, 4).Formula = "=""1.001"""
$ws.Cells.Item(12, 4).Copy()
$ws.Cells.Item(12, 4).PasteSpecial(-4163)
$ws.Cells.Item(12, 5).Formula = "=""  +0.34%  """
$ws.Cells.Item(12, 5).Copy()
$ws.Cells.Item(12, 5).PasteSpecial(-4163)
$ws.Cells.Item(13, 4).Formula = "=""21.20"""
$ws.Cells.Item(13, 4).Copy()
$ws.Cells.Item(13, 4).PasteSpecial(-4163)
$ws.Cells.Item(13, 5).Formula = "=""  +6.52%  """
$ws.Cells.Item(13, 5).Copy()
$ws.Cells.Item(13, 5).PasteSpecial(-4163)
$ws.Cells.Item(14, 4).Formula = "=""6.259"""
$ws.Cells.Item(14, 4).Copy()
$ws.Cells.Item(14, 4).PasteSpecial(-4163)
$ws.Cells.Item(14, 5).Formula = "=""  +6.55%  """
$ws.Cells.Item(14, 5).Copy()
$ws.Cells.Item(14, 5).PasteSpecial(-4163)
$ws.Cells.Item(15, 4).Formula = "=""1.811.50"""
$ws.Cells.Item(15, 4).Copy()
$ws.Cells.Item(15, 4).PasteSpecial(-4163)
$ws.Cells.Item(15, 5).Formula = "=""  +5.41%  """
$ws.Cells.Item(15, 5).Copy()
$ws.Cells.Item(15, 5).PasteSpecial(-4163)
$ws.Cells.Item(16, 4).Formula = "=""7.267"""
$ws.Cells.Item(16, 4).Copy()
$ws.Cells.Item(16, 4).PasteSpecial(-4163)
$ws.Cells.Item(16, 5).Formula = "=""  +6.64%  """
$ws.Cells.Item(16, 5).Copy()
$ws.Cells.Item(16, 5).PasteSpecial(-4163)
$ws.Cells.Item(17, 4).Formula = "=""92.00"""
$ws.Cells.Item(17, 4).Copy()
$ws.Cells.Item(17, 4).PasteSpecial(-4163)
$ws.Cells.Item(17, 5).Formula = "=""  +6.09%  """
$ws.Cells.Item(17, 5).Copy()
$ws.Cells.Item(17, 5).PasteSpecial(-4163)
$ws.Cells.Item(18, 4).Formula = "=""0.00001075"""
$ws.Cells.Item(18, 4).Copy()
$ws.Cells.Item(18, 4).PasteSpecial(-4163)
$ws.Cells.Item(18, 5).Formula = "=""  +3.75%  """
$ws.Cells.Item(18, 5).Copy()
$ws.Cells.Item(18, 5).PasteSpecial(-4163)
$ws.Cells.Item(19, 4).Formula = "=""0.06481"""
$ws.Cells.Item(19, 4).Copy()
$ws.Cells.Item(19, 4).PasteSpecial(-4163)
$ws.Cells.Item(19, 5).Formula = "=""  +1.38%  """
$ws.Cells.Item(19, 5).Copy()
$ws.Cells.Item(19, 5).PasteSpecial(-4163)
$ws.Cells.Item(20, 4).Formula = "=""1.001"""
$ws.Cells.Item(20, 4).Copy()
$ws.Cells.Item(20, 4).PasteSpecial(-4163)
$ws.Cells.Item(20, 5).Formula = "=""  +0.22%  """
$ws.Cells.Item(20, 5).Copy()
$ws.Cells.Item(20, 5).PasteSpecial(-4163)
$ws.Cells.Item(21, 4).Formula = "=""17.32"""
$ws.Cells.Item(21, 4).Copy()
$ws.Cells.Item(21, 4).PasteSpecial(-4163)
$ws.Cells.Item(21, 5).Formula = "=""  +4.73%  """
$ws.Cells.Item(21, 5).Copy()
$ws.Cells.Item(21, 5).PasteSpecial(-4163)
$ws.Cells.Item(22, 4).Formula = "=""6.010"""
$ws.Cells.Item(22, 4).Copy()
$ws.Cells.Item(22, 4).PasteSpecial(-4163)
$ws.Cells.Item(22, 5).Formula = "=""  +5.07%  """
$ws.Cells.Item(22, 5).Copy()
$ws.Cells.Item(22, 5).PasteSpecial(-4163)
$ws.Cells.Item(23, 4).Formula = "=""28.369.68"""
$ws.Cells.Item(23, 4).Copy()
$ws.Cells.Item(23, 4).PasteSpecial(-4163)
$ws.Cells.Item(23, 5).Formula = "=""  +5.18%  """
$ws.Cells.Item(23, 5).Copy()
$ws.Cells.Item(23, 5).PasteSpecial(-4163)
$ws.Cells.Item(24, 4).Formula = "=""11.32"""
$ws.Cells.Item(24, 4).Copy()
$ws.Cells.Item(24, 4).PasteSpecial(-4163)
$ws.Cells.Item(24, 5).Formula = "=""  +3.38%  """
$ws.Cells.Item(24, 5).Copy()
$ws.Cells.Item(24, 5).PasteSpecial(-4163)
$ws.Cells.Item(25, 4).Formula = "=""2.142"""
$ws.Cells.Item(25, 4).Copy()
$ws.Cells.Item(25, 4).PasteSpecial(-4163)
$ws.Cells.Item(25, 5).Formula = "=""  +4.17%  """
$ws.Cells.Item(25, 5).Copy()
$ws.Cells.Item(25, 5).PasteSpecial(-4163)
$ws.Cells.Item(26, 4).Formula = "=""20.83"""
$ws.Cells.Item(26, 4).Copy()
$ws.Cells.Item(26, 4).PasteSpecial(-4163)
$ws.Cells.Item(26, 5).Formula = "=""  +4.82%  """
$ws.Cells.Item(26, 5).Copy()
$ws.Cells.Item(26, 5).PasteSpecial(-4163)
$ws.Cells.Item(27, 4).Formula = "=""157.92"""
$ws.Cells.Item(27, 4).Copy()
$ws.Cells.Item(27, 4).PasteSpecial(-4163)
$ws.Cells.Item(27, 5).Formula = "=""  +2.35%  """
$ws.Cells.Item(27, 5).Copy()
$ws.Cells.Item(27, 5).PasteSpecial(-4163)
$ws.Cells.Item(28, 4).Formula = "=""2.450"""
$ws.Cells.Item(28, 4).Copy()
$ws.Cells.Item(28, 4).PasteSpecial(-4163)
$ws.Cells.Item(28, 5).Formula = "=""  +18.23%  """
$ws.Cells.Item(28, 5).Copy()
$ws.Cells.Item(28, 5).PasteSpecial(-4163)
$ws.Cells.Item(29, 4).Formula = "=""2.019.48"""
$ws.Cells.Item(29, 4).Copy()
$ws.Cells.Item(29, 4).PasteSpecial(-4163)
$ws.Cells.Item(29, 5).Formula = "=""  +5.40%  """
$ws.Cells.Item(29, 5).Copy()
$ws.Cells.Item(29, 5).PasteSpecial(-4163)
$ws.Cells.Item(30, 4).Formula = "=""124.07"""
$ws.Cells.Item(30, 4).Copy()
$ws.Cells.Item(30, 4).PasteSpecial(-4163)
$ws.Cells.Item(30, 5).Formula = "=""  +3.33%  """
$ws.Cells.Item(30, 5).Copy()
$ws.Cells.Item(30, 5).PasteSpecial(-4163)
$ws.Cells.Item(31, 4).Formula = "=""1.160"""
$ws.Cells.Item(31, 4).Copy()
$ws.Cells.Item(31, 4).PasteSpecial(-4163)
$ws.Cells.Item(31, 5).Formula = "=""  +10.87%  """
$ws.Cells.Item(31, 5).Copy()
$ws.Cells.Item(31, 5).PasteSpecial(-4163)
$ws.Cells.Item(32, 4).Formula = "=""0.1065"""
$ws.Cells.Item(32, 4).Copy()
$ws.Cells.Item(32, 4).PasteSpecial(-4163)
$ws.Cells.Item(32, 5).Formula = "=""  +14.61%  """
$ws.Cells.Item(32, 5).Copy()
$ws.Cells.Item(32, 5).PasteSpecial(-4163)
$ws.Cells.Item(33, 4).Formula = "=""5.789"""
$ws.Cells.Item(33, 4).Copy()
$ws.Cells.Item(33, 4).PasteSpecial(-4163)
$ws.Cells.Item(33, 5).Formula = "=""  +7.74%  """
$ws.Cells.Item(33, 5).Copy()
$ws.Cells.Item(33, 5).PasteSpecial(-4163)
$ws.Cells.Item(34, 4).Formula = "=""3.629"""
$ws.Cells.Item(34, 4).Copy()
$ws.Cells.Item(34, 4).PasteSpecial(-4163)
$ws.Cells.Item(34, 5).Formula = "=""  +1.43%  """
$ws.Cells.Item(34, 5).Copy()
$ws.Cells.Item(34, 5).PasteSpecial(-4163)
$ws.Cells.Item(35, 4).Formula = "=""0.2202"""
$ws.Cells.Item(35, 4).Copy()
$ws.Cells.Item(35, 4).PasteSpecial(-4163)
$ws.Cells.Item(35, 5).Formula = "=""  +10.82%  """
$ws.Cells.Item(35, 5).Copy()
$ws.Cells.Item(35, 5).PasteSpecial(-4163)
$ws.Cells.Item(36, 4).Formula = "=""8.919"""
$ws.Cells.Item(36, 4).Copy()
$ws.Cells.Item(36, 4).PasteSpecial(-4163)
$ws.Cells.Item(36, 5).Formula = "=""  +20.15%  """
$ws.Cells.Item(36, 5).Copy()
$ws.Cells.Item(36, 5).PasteSpecial(-4163)
$ws.Cells.Item(37, 4).Formula = "=""0.02323"""
$ws.Cells.Item(37, 4).Copy()
$ws.Cells.Item(37, 4).PasteSpecial(-4163)
$ws.Cells.Item(37, 5).Formula = "=""  +6.49%  """
$ws.Cells.Item(37, 5).Copy()
$ws.Cells.Item(37, 5).PasteSpecial(-4163)
$ws.Cells.Item(38, 4).Formula = "=""11.69"""
$ws.Cells.Item(38, 4).Copy()
$ws.Cells.Item(38, 4).PasteSpecial(-4163)
$ws.Cells.Item(38, 5).Formula = "=""  +6.71%  """
$ws.Cells.Item(38, 5).Copy()
$ws.Cells.Item(38, 5).PasteSpecial(-4163)
$ws.Cells.Item(39, 4).Formula = "=""0.06118"""
$ws.Cells.Item(39, 4).Copy()
$ws.Cells.Item(39, 4).PasteSpecial(-4163)
$ws.Cells.Item(39, 5).Formula = "=""  +3.76%  """
$ws.Cells.Item(39, 5).Copy()
$ws.Cells.Item(39, 5).PasteSpecial(-4163)
$ws.Cells.Item(40, 4).Formula = "=""5.048"""
$ws.Cells.Item(40, 4).Copy()
$ws.Cells.Item(40, 4).PasteSpecial(-4163)
$ws.Cells.Item(40, 5).Formula = "=""  +6.29%  """
$ws.Cells.Item(40, 5).Copy()
$ws.Cells.Item(40, 5).PasteSpecial(-4163)
$ws.Cells.Item(41, 4).Formula = "=""0.6410"""
$ws.Cells.Item(41, 4).Copy()
$ws.Cells.Item(41, 4).PasteSpecial(-4163)
$ws.Cells.Item(41, 5).Formula = "=""  +7.10%  """
$ws.Cells.Item(41, 5).Copy()
$ws.Cells.Item(41, 5).PasteSpecial(-4163)
$ws.Cells.Item(42, 4).Formula = "=""1.164"""
$ws.Cells.Item(42, 4).Copy()
$ws.Cells.Item(42, 4).PasteSpecial(-4163)
$ws.Cells.Item(42, 5).Formula = "=""  +4.47%  """
$ws.Cells.Item(42, 5).Copy()
$ws.Cells.Item(42, 5).PasteSpecial(-4163)
$ws.Cells.Item(43, 4).Formula = "=""1.000"""
$ws.Cells.Item(43, 4).Copy()
$ws.Cells.Item(43, 4).PasteSpecial(-4163)
$ws.Cells.Item(43, 5).Formula = "=""  +0.18%  """
$ws.Cells.Item(43, 5).Copy()
$ws.Cells.Item(43, 5).PasteSpecial(-4163)
$ws.Cells.Item(44, 4).Formula = "=""1.379"""
$ws.Cells.Item(44, 4).Copy()
$ws.Cells.Item(44, 4).PasteSpecial(-4163)
$ws.Cells.Item(44, 5).Formula = "=""  -3.22%  """
$ws.Cells.Item(44, 5).Copy()
$ws.Cells.Item(44, 5).PasteSpecial(-4163)
$ws.Cells.Item(45, 4).Formula = "=""13.43"""
$ws.Cells.Item(45, 4).Copy()
$ws.Cells.Item(45, 4).PasteSpecial(-4163)
$ws.Cells.Item(45, 5).Formula = "=""  +4.98%  """
$ws.Cells.Item(45, 5).Copy()
$ws.Cells.Item(45, 5).PasteSpecial(-4163)
$ws.Cells.Item(46, 4).Formula = "=""0.6008"""
$ws.Cells.Item(46, 4).Copy()
$ws.Cells.Item(46, 4).PasteSpecial(-4163)
$ws.Cells.Item(46, 5).Formula = "=""  +6.98%  """
$ws.Cells.Item(46, 5).Copy()
$ws.Cells.Item(46, 5).PasteSpecial(-4163)
$ws.Cells.Item(47, 4).Formula = "=""3.705"""
$ws.Cells.Item(47, 4).Copy()
$ws.Cells.Item(47, 4).PasteSpecial(-4163)
$ws.Cells.Item(47, 5).Formula = "=""  +3.65%  """
$ws.Cells.Item(47, 5).Copy()
$ws.Cells.Item(47, 5).PasteSpecial(-4163)
$ws.Cells.Item(48, 4).Formula = "=""122.11"""
$ws.Cells.Item(48, 4).Copy()
$ws.Cells.Item(48, 4).PasteSpecial(-4163)
$ws.Cells.Item(48, 5).Formula = "=""  +2.19%  """
$ws.Cells.Item(48, 5).Copy()
$ws.Cells.Item(48, 5).PasteSpecial(-4163)
$ws.Cells.Item(49, 4).Formula = "=""1.946"""
$ws.Cells.Item(49, 4).Copy()
$ws.Cells.Item(49, 4).PasteSpecial(-4163)
$ws.Cells.Item(49, 5).Formula = "=""  +5.80%  """
$ws.Cells.Item(49, 5).Copy()
$ws.Cells.Item(49, 5).PasteSpecial(-4163)
$ws.Cells.Item(50, 4).Formula = "=""1.150"""
$ws.Cells.Item(50, 4).Copy()
$ws.Cells.Item(50, 4).PasteSpecial(-4163)
$ws.Cells.Item(50, 5).Formula = "=""  +5.09%  """
$ws.Cells.Item(50, 5).Copy()
$ws.Cells.Item(50, 5).PasteSpecial(-4163)
$ws.Cells.Item(51, 4).Formula = "=""0.06868"""
$ws.Cells.Item(51, 4).Copy()
$ws.Cells.Item(51, 4).PasteSpecial(-4163)
$ws.Cells.Item(51, 5).Formula = "=""  +3.35%  """
$ws.Cells.Item(51, 5).Copy()
$ws.Cells.Item(51, 5).PasteSpecial(-4163)
$excel.CutCopyMode = $false
